$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire third row (it is removed from the data set)
$ws.Rows.Item(3).Delete()

# Update row 2 values with the new data set (order matches author's entry order)
$ws.Range("A2").Value = 11
$ws.Range("G2").Value = "Test1"
$ws.Range("I2").Value = "ABC"
$ws.Range("C2").Value = "123456"
$ws.Range("F2").Value = "10.17.18.88"
$ws.Range("E2").Value = "Setup-1"
$ws.Range("D2").Value = "Ast1"
$ws.Range("H2").Value = "7"
$ws.Range("B2").Value = "1"

# Update selection to match the new active cell
$ws.Range("H2").Select()
